$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Le Doppio Impasto" block (rows 3-6): replace old pizza entries with the
# new "Variazioni" style entries (name + adjustment only, description cleared)
$ws.Range("B3").Value = "Pizza Baby"
$ws.Range("C3").ClearContents()
$ws.Range("D3").Value = -1

$ws.Range("B4").Value = "Pizza Doppio Impasto"
$ws.Range("C4").ClearContents()
$ws.Range("D4").Value = 1.5

$ws.Range("B5").Value = "Pizza con farina al Farro Integrale"
$ws.Range("C5").ClearContents()
$ws.Range("D5").Value = 1

$ws.Range("B6").Value = "Pizza con farina Grano Khorasan"
$ws.Range("C6").ClearContents()
$ws.Range("D6").Value = 1

# Rows 7-10 no longer have entries - clear them out entirely
$ws.Range("B7:D7").ClearContents()
$ws.Range("B8:D8").ClearContents()
$ws.Range("B9:D9").ClearContents()
$ws.Range("B10:D10").ClearContents()

# Scroll position / selection moved
$ws.Range("K9").Select()
